$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.337.56"
$ws.Range("E2").Value = "  -2.19%  "
$ws.Range("D3").Value = "3.486.99"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.85%  "
$ws.Range("D7").Value = "3.487.08"
$ws.Range("E7").Value = "  -3.13%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "4.078.19"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("E15").Value = "  -4.65%  "
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "3.484.27"
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("D18").Value = "64.379.22"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.576"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").Value = "3.625.74"
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -9.08%  "
$ws.Range("E28").Value = "  -7.15%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.72%  "
$ws.Range("E31").Value = "  -6.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").Value = "3.486.01"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "170.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.811"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.43%  "
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.436.28"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0267"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.37%  "
